$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.238266825675964
$ws.Range("B1").Value = 2.304747343063354
$ws.Range("D1").Value = 1.41278600692749
$ws.Range("E1").Value = 0.8835155367851257
